$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update loading_percent values for rows 2-25 (columns B, C, D, E, F, M)

# Row 2
$ws.Range("B2").Value = 14.61399645465955
$ws.Range("C2").Value = 7.967680387001399
$ws.Range("D2").Value = 7.919997271376043
$ws.Range("E2").Value = 10.35691967998542
$ws.Range("F2").Value = 45.35077606733255
$ws.Range("M2").Value = 15.76558620686709

# Row 3
$ws.Range("B3").Value = 14.20314885177125
$ws.Range("C3").Value = 7.550689470232961
$ws.Range("D3").Value = 7.757603831157018
$ws.Range("E3").Value = 10.192126637831
$ws.Range("F3").Value = 43.91823071386094
$ws.Range("M3").Value = 15.56563962829293

# Row 4
$ws.Range("B4").Value = 13.95270309439362
$ws.Range("C4").Value = 7.28749690518789
$ws.Range("D4").Value = 7.656282840507675
$ws.Range("E4").Value = 10.09202247293471
$ws.Range("F4").Value = 43.02161795793959
$ws.Range("M4").Value = 15.44798760829665

# Row 5
$ws.Range("B5").Value = 13.85130646400723
$ws.Range("C5").Value = 7.178642268385903
$ws.Range("D5").Value = 7.614621340415807
$ws.Range("E5").Value = 10.0515433978895
$ws.Range("F5").Value = 42.65244284402643
$ws.Range("M5").Value = 15.40138909386624

# Row 6
$ws.Range("B6").Value = 13.83451573891822
$ws.Range("C6").Value = 7.160476716383031
$ws.Range("D6").Value = 7.607681964732704
$ws.Range("E6").Value = 10.04484209197425
$ws.Range("F6").Value = 42.5909269977293
$ws.Range("M6").Value = 15.39373438137247

# Row 7
$ws.Range("B7").Value = 13.95133266733582
$ws.Range("C7").Value = 7.286035053115782
$ws.Range("D7").Value = 7.655722442757293
$ws.Range("E7").Value = 10.0914752297896
$ws.Range("F7").Value = 43.01665384794129
$ws.Range("M7").Value = 15.44735364028418

# Row 8
$ws.Range("B8").Value = 14.47208025210429
$ws.Range("C8").Value = 7.825497506956293
$ws.Range("D8").Value = 7.864357575617871
$ws.Range("E8").Value = 10.29990080036252
$ws.Range("F8").Value = 44.86064629784929
$ws.Range("M8").Value = 15.69561931726194

# Row 9
$ws.Range("B9").Value = 15.49905265961384
$ws.Range("C9").Value = 8.835613099371578
$ws.Range("D9").Value = 8.259337938602014
$ws.Range("E9").Value = 10.71529438401304
$ws.Range("F9").Value = 48.32227871434101
$ws.Range("M9").Value = 16.2203090866832

# Row 10
$ws.Range("B10").Value = 16.2462653797859
$ws.Range("C10").Value = 9.649016824642068
$ws.Range("D10").Value = 8.539166244345269
$ws.Range("E10").Value = 11.02202505737253
$ws.Range("F10").Value = 50.74812194904771
$ws.Range("M10").Value = 16.6250232440161

# Row 11
$ws.Range("B11").Value = 16.58255383145425
$ws.Range("C11").Value = 9.998869313032468
$ws.Range("D11").Value = 8.663879833171404
$ws.Range("E11").Value = 11.16137461780933
$ws.Range("F11").Value = 51.82202434326884
$ws.Range("M11").Value = 16.81246874501927

# Row 12
$ws.Range("B12").Value = 16.70922376080208
$ws.Range("C12").Value = 10.1284625450855
$ws.Range("D12").Value = 8.710708272926651
$ws.Range("E12").Value = 11.21407371468427
$ws.Range("F12").Value = 52.22412176610551
$ws.Range("M12").Value = 16.88385865377017

# Row 13
$ws.Range("B13").Value = 16.68197546769235
$ws.Range("C13").Value = 10.10068057884691
$ws.Range("D13").Value = 8.700641032649143
$ws.Range("E13").Value = 11.2027277864032
$ws.Range("F13").Value = 52.13773048566733
$ws.Range("M13").Value = 16.8684665334769

# Row 14
$ws.Range("B14").Value = 16.59298941054043
$ws.Range("C14").Value = 10.00958883940207
$ws.Range("D14").Value = 8.667740557828131
$ws.Range("E14").Value = 11.16571184023085
$ws.Range("F14").Value = 51.85519828490906
$ws.Range("M14").Value = 16.81833426045351

# Row 15
$ws.Range("B15").Value = 16.5383906011999
$ws.Range("C15").Value = 9.953416751492657
$ws.Range("D15").Value = 8.647535488905582
$ws.Range("E15").Value = 11.14302819007449
$ws.Range("F15").Value = 51.68153600634559
$ws.Range("M15").Value = 16.7876778568237

# Row 16
$ws.Range("B16").Value = 16.22420174549663
$ws.Range("C16").Value = 9.625748283753479
$ws.Range("D16").Value = 8.530961650494485
$ws.Range("E16").Value = 11.01291081942567
$ws.Range("F16").Value = 50.67731682762388
$ws.Range("M16").Value = 16.61283428138539

# Row 17
$ws.Range("B17").Value = 16.03041517521312
$ws.Range("C17").Value = 9.419576630585794
$ws.Range("D17").Value = 8.458766635722899
$ws.Range("E17").Value = 10.93300984206414
$ws.Range("F17").Value = 50.05344450202808
$ws.Range("M17").Value = 16.50637739833046

# Row 18
$ws.Range("B18").Value = 15.91862167939043
$ws.Range("C18").Value = 9.299092943341565
$ws.Range("D18").Value = 8.417000454920617
$ws.Range("E18").Value = 10.88703811801311
$ws.Range("F18").Value = 49.69183718955457
$ws.Range("M18").Value = 16.44546478204059

# Row 19
$ws.Range("B19").Value = 15.88071812535194
$ws.Range("C19").Value = 9.257972703874715
$ws.Range("D19").Value = 8.402818502962269
$ws.Range("E19").Value = 10.87147167780883
$ws.Range("F19").Value = 49.56893687434836
$ws.Range("M19").Value = 16.42489772583077

# Row 20
$ws.Range("B20").Value = 16.0510796689118
$ws.Range("C20").Value = 9.441720432928031
$ws.Range("D20").Value = 8.466477112935912
$ws.Range("E20").Value = 10.94151728562525
$ws.Range("F20").Value = 50.12014603780022
$ws.Range("M20").Value = 16.51767746645624

# Row 21
$ws.Range("B21").Value = 16.6191462599542
$ws.Range("C21").Value = 10.03642297790318
$ws.Range("D21").Value = 8.677415227164422
$ws.Range("E21").Value = 11.17658654534433
$ws.Range("F21").Value = 51.93831103557767
$ws.Range("M21").Value = 16.83304881029606

# Row 22
$ws.Range("B22").Value = 16.98641298227503
$ws.Range("C22").Value = 10.40826654424318
$ws.Range("D22").Value = 8.812945143224843
$ws.Range("E22").Value = 11.32979350999234
$ws.Range("F22").Value = 53.09985357738969
$ws.Range("M22").Value = 17.04150974892217

# Row 23
$ws.Range("B23").Value = 16.79080907584903
$ws.Range("C23").Value = 10.21134222308647
$ws.Range("D23").Value = 8.740831855208393
$ws.Range("E23").Value = 11.24807677111075
$ws.Range("F23").Value = 52.48245505621434
$ws.Range("M23").Value = 16.93005894861458

# Row 24
$ws.Range("B24").Value = 16.0417384324588
$ws.Range("C24").Value = 9.431715303496803
$ws.Range("D24").Value = 8.462992014939777
$ws.Range("E24").Value = 10.93767117795147
$ws.Range("F24").Value = 50.08999939115579
$ws.Range("M24").Value = 16.51256779422907

# Row 25
$ws.Range("B25").Value = 15.22177272251135
$ws.Range("C25").Value = 8.558233268203578
$ws.Range("D25").Value = 8.154185318953635
$ws.Range("E25").Value = 10.60246756448853
$ws.Range("F25").Value = 47.40489543236756
$ws.Range("M25").Value = 16.07470601068318
